$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.26"
$ws.Range("E2").Value = "'-0.78%"
$ws.Range("D3").Value = "'27.33"
$ws.Range("E3").Value = "'4.07%"
$ws.Range("D4").Value = "'5.113"
$ws.Range("E4").Value = "'0.63%"
$ws.Range("D5").Value = "'0.05685"
$ws.Range("E5").Value = "'1.50%"
$ws.Range("D6").Value = "'6.519"
$ws.Range("E6").Value = "'0.69%"
$ws.Range("E7").Value = "'0.88%"
$ws.Range("E8").Value = "'1.94%"
$ws.Range("B9").Value = "'WazirX"
$ws.Range("C9").Value = "'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").Value = "'0.1332"
$ws.Range("E9").Value = "'-0.53%"
$ws.Range("D10").Value = "'0.06955"
$ws.Range("E10").Value = "'-0.42%"
$ws.Range("D11").Value = "'0.02857"
$ws.Range("E11").Value = "'1.70%"
$ws.Range("D12").Value = "'0.09386"
$ws.Range("E12").Value = "'0.00%"
$ws.Range("D13").Value = "'0.001518"
$ws.Range("E13").Value = "'0.34%"
$ws.Range("D14").Value = "'0.04071"
$ws.Range("E14").Value = "'-12.88%"
$ws.Range("B15").Value = "'One"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15").Value = "'0.0005979"
$ws.Range("E15").Value = "'-0.50%"
$ws.Range("B16").Value = "'TigerCash"
$ws.Range("C16").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006214"
$ws.Range("E16").Value = "'1.07%"
$ws.Range("B17").Value = "'LEO"
$ws.Range("C17").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.512"
$ws.Range("E17").Value = "'-2.70%"
$ws.Range("B18").Value = "'GateToken"
$ws.Range("C18").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'3.009"
$ws.Range("E18").Value = "'-0.32%"
$ws.Range("B19").Value = "'BTSEToken"
$ws.Range("C19").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.317"
$ws.Range("E19").Value = "'12.74%"
$ws.Range("B20").Value = "'BitpandaEcosystemToken"
$ws.Range("C20").Value = "'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3164"
$ws.Range("E20").Value = "'1.52%"
$ws.Range("D21").Value = "'0.03222"
$ws.Range("E21").Value = "'0.57%"
$ws.Range("D23").Value = "'3.552"
$ws.Range("E23").Value = "'-5.20%"
$ws.Range("E24").Value = "'1.78%"
$ws.Range("D25").Value = "'0.001216"
$ws.Range("E25").Value = "'-2.16%"
$ws.Range("D26").Value = "'0.004466"
$ws.Range("E27").Value = "'22.92%"
$ws.Range("D28").Value = "'0.0001406"
$ws.Range("E28").Value = "'-27.44%"
$ws.Range("E40").Value = "'1.66%"
$ws.Range("D41").Value = "'0.005920"
$ws.Range("E41").Value = "'-3.42%"
$ws.Range("E42").Value = "'0.18%"
$ws.Range("D43").Value = "'0.001800"
$ws.Range("E43").Value = "'-28.00%"
$ws.Range("D44").Value = "'0.009714"
$ws.Range("E44").Value = "'17.16%"
$ws.Range("D45").Value = "'0.00005101"
$ws.Range("E45").Value = "'-5.03%"
$ws.Range("E46").Value = "'0.00%"
$ws.Range("E47").Value = "'-30.36%"
$ws.Range("E48").Value = "'-2.80%"
$ws.Range("E49").Value = "'0.00%"
$ws.Range("E50").Value = "'0.00%"
